$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder player rows to match the updated roster layout.
# Row 3 <- old row 14 (Jalen Green)
$ws.Range("A3").Value = "Jalen Green"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Houston Rockets"

# Row 4 <- old row 15 (Jalen Suggs)
$ws.Range("A4").Value = "Jalen Suggs"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Orlando Magic"

# Row 5 <- old row 4 (Chris Paul)
$ws.Range("A5").Value = "Chris Paul"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "San Antonio Spurs"

# Row 11 <- old row 16 (Clint Capela)
$ws.Range("A11").Value = "Clint Capela"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Atlanta Hawks"

# Row 14 <- old row 3 (Klay Thompson)
$ws.Range("A14").Value = "Klay Thompson"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Dallas Mavericks"

# Row 15 <- old row 5 (Jaylen Brown)
$ws.Range("A15").Value = "Jaylen Brown"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Boston Celtics"

# Row 16 <- old row 11 (Jakob Poeltl)
$ws.Range("A16").Value = "Jakob Poeltl"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Toronto Raptors"
